# Edit script: append new rows to three_line, two_line, and ph_pl_breakout_line sheets
# Column layout reference:
#   three_line (sheet1):          A detected_date(date) B breakout_date C Time_Frame D stockname
#                                  E date1(date) F value1 G date2(date) H value2 I date3 J value3
#                                  K buyORsell L "Date Time"
#   two_line (sheet2):             A detected_date(date) B breakout_date C Time_Frame D stockname
#                                  E date1(date) F value1 G date2(date) H value2 I buyORsell J "Date Time"
#   ph_pl_breakout_line (sheet3):  A stockname B Datetime(date) C High D Low E Close F isPivot
#                                  G PHorPLValue H time_frame I TdyDate J TdyClose K PClose L "Date Time"

$wb = $excel.ActiveWorkbook


function Add-Rows {
    param($ws, $rows, $dateCols, $numCols, $dateFormat)

    foreach ($r in $rows) {
        $rowNum = $r[0]
        for ($ci = 0; $ci -lt ($r.Count - 1); $ci++) {
            $col = $ci + 1
            $raw = $r[$ci + 1]
            $cell = $ws.Cells.Item($rowNum, $col)
            if ($dateCols -contains $col) {
                $cell.Value = [double]$raw
                $cell.NumberFormat = $dateFormat
            } elseif ($numCols -contains $col) {
                $cell.Value = [double]$raw
            } else {
                $cell.Value = [string]$raw
            }
        }
    }
}

# ---- Sheet: three_line ----
$wsThree = $wb.Worksheets.Item("three_line")
$threeDateFormat = $wsThree.Range("A120").NumberFormat

$sheet1Rows = @(
    ,@(121, "45383", "10-06-2024 00:00:00", "week", "JSWSTEEL.NS", "45285", "895.75", "45341", "850", "0", "0", "High", "10/06/2024 14:03:06")
    ,@(122, "45355", "10-06-2024 00:00:00", "week", "TECHM.NS", "45187", "1320", "45313", "1416.300048828125", "0", "0", "High", "10/06/2024 14:03:06")
    ,@(123, "45348", "10-06-2024 00:00:00", "week", "PIDILITIND.NS", "44816", "2918.949951171875", "45306", "2805", "0", "0", "High", "10/06/2024 14:03:06")
    ,@(124, "45432", "10-06-2024 00:00:00", "week", "IRB.NS", "45327", "72", "45390", "73.09999847412109", "0", "0", "High", "10/06/2024 14:03:06")
    ,@(125, "45369", "10-06-2024 00:00:00", "week", "JISLJALEQS.NS", "45173", "69.80000305175781", "45250", "73.44999694824219", "0", "0", "High", "10/06/2024 14:03:06")
    ,@(126, "44319", "10-06-2024 00:00:00", "week", "BLUECLOUDS.BO", "43549", "13.10999965667725", "44200", "12.52000045776367", "0", "0", "High", "10/06/2024 14:03:06")
)

Add-Rows $wsThree $sheet1Rows @(1,5,7) @(6,8,9,10) $threeDateFormat

# ---- Sheet: two_line ----
$wsTwo = $wb.Worksheets.Item("two_line")
$twoDateFormat = $wsTwo.Range("A27").NumberFormat

$sheet2Rows = @(
    ,@(28, "45292", "10-06-2024 00:00:00", "week", "BAJAJELEC.NS", "45026", "974.8753051757812", "45250", "975", "Low", "10/06/2024 14:03:06")
)

Add-Rows $wsTwo $sheet2Rows @(1,5,7) @(6,8) $twoDateFormat

# ---- Sheet: ph_pl_breakout_line ----
$wsPh = $wb.Worksheets.Item("ph_pl_breakout_line")
$phDateFormat = $wsPh.Range("B502").NumberFormat

$sheet3Rows = @(
    ,@(503, "CHOLAFIN.NS", "45411", "1352.599975585938", "1161.150024414062", "1309.699951171875", "High", "1352.599975585938", "week", "10-06-2024 00:00:00", "1369", "1351", "10/06/2024 14:03:06")
    ,@(504, "ULTRACEMCO.NS", "45285", "10526", "9969", "10503.0498046875", "High", "10526", "week", "10-06-2024 00:00:00", "10907.9501953125", "10519.7998046875", "10/06/2024 14:03:06")
    ,@(505, "CIPLA.NS", "45362", "1519", "1449", "1488.050048828125", "High", "1519", "week", "10-06-2024 00:00:00", "1540.849975585938", "1514", "10/06/2024 14:03:06")
    ,@(506, "SHREECEM.NS", "45215", "27298.94921875", "25900", "26084.650390625", "High", "27298.94921875", "week", "10-06-2024 00:00:00", "27450", "26250", "10/06/2024 14:03:06")
    ,@(507, "SHREECEM.NS", "45383", "26743.30078125", "25699", "25811.05078125", "High", "26743.30078125", "week", "10-06-2024 00:00:00", "27450", "26250", "10/06/2024 14:03:06")
    ,@(508, "PGHH.NS", "45376", "17050.900390625", "16111", "16928.44921875", "High", "17050.900390625", "week", "10-06-2024 00:00:00", "17064.94921875", "16947.19921875", "10/06/2024 14:03:06")
    ,@(509, "VBL.NS", "45348", "1561.949951171875", "1345", "1417.849975585938", "High", "1561.949951171875", "week", "10-06-2024 00:00:00", "1593.699951171875", "1528", "10/06/2024 14:03:06")
    ,@(510, "BAJAJCON.NS", "45222", "256", "226.4499969482422", "240", "High", "256", "week", "10-06-2024 00:00:00", "260.3999938964844", "255.8999938964844", "10/06/2024 14:03:06")
    ,@(511, "MAHLOG.NS", "45299", "460", "435", "437.9500122070312", "High", "460", "week", "10-06-2024 00:00:00", "473.8999938964844", "458.25", "10/06/2024 14:03:06")
    ,@(512, "APLAPOLLO.NS", "45362", "1660", "1440", "1474.699951171875", "High", "1660", "week", "10-06-2024 00:00:00", "1667.5", "1631.5", "10/06/2024 14:03:06")
    ,@(513, "JISLJALEQS.NS", "45173", "69.80000305175781", "61.59999847412109", "65.90000152587891", "High", "69.80000305175781", "week", "10-06-2024 00:00:00", "76.34999847412109", "67.69999694824219", "10/06/2024 14:03:06")
    ,@(514, "JISLJALEQS.NS", "45250", "73.44999694824219", "67.05000305175781", "71.19999694824219", "High", "73.44999694824219", "week", "10-06-2024 00:00:00", "76.34999847412109", "67.69999694824219", "10/06/2024 14:03:06")
    ,@(515, "JISLJALEQS.NS", "45327", "68.59999847412109", "61.65000152587891", "66.05000305175781", "High", "68.59999847412109", "week", "10-06-2024 00:00:00", "76.34999847412109", "67.69999694824219", "10/06/2024 14:03:06")
    ,@(516, "CHEMPLASTS.NS", "45278", "527.7999877929688", "472.3999938964844", "513", "High", "527.7999877929688", "week", "10-06-2024 00:00:00", "559", "495", "10/06/2024 14:03:06")
    ,@(517, "CHEMPLASTS.NS", "45348", "522.2000122070312", "464", "491.9500122070312", "High", "522.2000122070312", "week", "10-06-2024 00:00:00", "559", "495", "10/06/2024 14:03:06")
    ,@(518, "CHEMPLASTS.NS", "45411", "527", "476.2000122070312", "479.3999938964844", "High", "527", "week", "10-06-2024 00:00:00", "559", "495", "10/06/2024 14:03:06")
    ,@(519, "CHEMPLASTS.NS", "45145", "464.9500122070312", "416.1499938964844", "444.0499877929688", "Low", "416.1499938964844", "week", "10-06-2024 00:00:00", "402.7999877929688", "441.8500061035156", "10/06/2024 14:03:06")
    ,@(520, "CHEMPLASTS.NS", "45250", "457.1000061035156", "425", "444.6000061035156", "Low", "425", "week", "10-06-2024 00:00:00", "402.7999877929688", "441.8500061035156", "10/06/2024 14:03:06")
    ,@(521, "CHEMPLASTS.NS", "45334", "480.9500122070312", "413.4500122070312", "450.75", "Low", "413.4500122070312", "week", "10-06-2024 00:00:00", "402.7999877929688", "441.8500061035156", "10/06/2024 14:03:06")
    ,@(522, "PREMEXPLN.NS", "45411", "2780", "2326.64990234375", "2431.25", "High", "2780", "week", "10-06-2024 00:00:00", "3045.449951171875", "2688", "10/06/2024 14:03:06")
    ,@(523, "TATACOMM.NS", "45278", "1842", "1688.050048828125", "1730.150024414062", "High", "1842", "week", "10-06-2024 00:00:00", "1877.599975585938", "1821.400024414062", "10/06/2024 14:03:06")
    ,@(524, "JKPAPER.NS", "45215", "415.2999877929688", "391", "400.25", "High", "415.2999877929688", "week", "10-06-2024 00:00:00", "429", "407.5499877929688", "10/06/2024 14:03:06")
    ,@(525, "HTMEDIA.NS", "45124", "26.29999923706055", "20.85000038146973", "24.35000038146973", "High", "26.29999923706055", "week", "10-06-2024 00:00:00", "26.75", "26.29999923706055", "10/06/2024 14:03:06")
)

Add-Rows $wsPh $sheet3Rows @(2) @(3,4,5,7,10,11) $phDateFormat

Write-Host "Rows appended successfully"
